$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.410.61'
$ws.Range("E2").Value = '  +3.37%  '
$ws.Range("D3").Value = '1.803.25'
$ws.Range("E3").Value = '  +4.71%  '
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +0.79%  '
$ws.Range("D5").Value = '333.66'
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").Value = '0.3793'
$ws.Range("E7").Value = '  +2.29%  '
$ws.Range("D8").Value = '0.3494'
$ws.Range("E8").Value = '  +4.31%  '
$ws.Range("D9").Value = '49.15'
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("D10").Value = '1.212'
$ws.Range("E10").Value = '  +2.63%  '
$ws.Range("D11").Value = '0.07598'
$ws.Range("E11").Value = '  +3.17%  '
$ws.Range("D12").Value = '1.009'
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("D13").Value = '21.94'
$ws.Range("E13").Value = '  +9.29%  '
$ws.Range("D14").Value = '6.519'
$ws.Range("E14").Value = '  +2.27%  '
$ws.Range("D15").Value = '1.809.19'
$ws.Range("E15").Value = '  +5.11%  '
$ws.Range("D16").Value = '7.094'
$ws.Range("E16").Value = '  +1.02%  '
$ws.Range("D17").Value = '0.00001105'
$ws.Range("E17").Value = '  +3.27%  '
$ws.Range("D18").Value = '0.06693'
$ws.Range("E18").Value = '  +1.02%  '
$ws.Range("D19").Value = '85.23'
$ws.Range("E19").Value = '  +3.50%  '
$ws.Range("D20").Value = '1.005'
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").Value = '17.37'
$ws.Range("E21").Value = '  +4.99%  '
$ws.Range("D22").Value = '6.469'
$ws.Range("E22").Value = '  +6.04%  '
$ws.Range("D23").Value = '27.447.11'
$ws.Range("E23").Value = '  +3.78%  '
$ws.Range("D24").Value = '12.63'
$ws.Range("E24").Value = '  -1.33%  '
$ws.Range("D25").Value = '2.445'
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("D26").Value = '2.594'
$ws.Range("E26").Value = '  +8.64%  '
$ws.Range("D27").Value = '21.67'
$ws.Range("E27").Value = '  +11.88%  '
$ws.Range("D28").Value = '1.450'
$ws.Range("E28").Value = '  +3.42%  '
$ws.Range("D29").Value = '150.27'
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("D30").Value = '2.013.68'
$ws.Range("E30").Value = '  +5.05%  '
$ws.Range("D31").Value = '134.35'
$ws.Range("E31").Value = '  +2.95%  '
$ws.Range("D32").Value = '4.082'
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("D33").Value = '6.116'
$ws.Range("E33").Value = '  +2.76%  '
$ws.Range("D34").Value = '0.08695'
$ws.Range("E34").Value = '  +1.57%  '
$ws.Range("D35").Value = '13.48'
$ws.Range("E35").Value = '  +6.55%  '
$ws.Range("D36").Value = '1.687'
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("D37").Value = '5.501'
$ws.Range("E37").Value = '  +2.84%  '
$ws.Range("D38").Value = '0.6865'
$ws.Range("E38").Value = '  +11.32%  '
$ws.Range("D39").Value = '8.980'
$ws.Range("E39").Value = '  +6.49%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.06398'
$ws.Range("E40").Value = '  +3.30%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.2220'
$ws.Range("E41").Value = '  +3.16%  '
$ws.Range("D42").Value = '0.02363'
$ws.Range("E42").Value = '  +1.85%  '
$ws.Range("D43").Value = '1.291'
$ws.Range("E43").Value = '  +5.50%  '
$ws.Range("D44").Value = '14.44'
$ws.Range("E44").Value = '  +3.10%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = '1.005'
$ws.Range("E45").Value = '  +0.39%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6428'
$ws.Range("E46").Value = '  +7.95%  '
$ws.Range("D47").Value = '3.845'
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("D48").Value = '2.137'
$ws.Range("E48").Value = '  +5.12%  '
$ws.Range("D49").Value = '131.37'
$ws.Range("E49").Value = '  +2.80%  '
$ws.Range("D50").Value = '0.07245'
$ws.Range("E50").Value = '  +1.02%  '
$ws.Range("D51").Value = '79.85'
$ws.Range("E51").Value = '  +4.34%  '
